# MASS_INTERVIEW_HISTORY_DATA.xlsx — "fixes in HTML report"
#
# Appends newly-run sprint/interview rows to the AMSIN, BETA and AMS
# history sheets, and fixes the styling of the previously-last AMS row
# (20 / "152_livetest") which had been missing the standard body-row
# style (and had a very slightly stale run-time timestamp) because it
# was the most recently appended row at the time it was written.

$wb = $excel.ActiveWorkbook

$wsAmsin = $wb.Worksheets.Item("AMSIN")
$wsBeta  = $wb.Worksheets.Item("BETA")
$wsAms   = $wb.Worksheets.Item("AMS")

# ---------------------------------------------------------------------
# AMS!A20:G20 ("152_livetest") — re-apply the normal body-row style
# (ClearContents resets the cell so the next write re-inherits the
# sheet's column default style) and refresh B20's run-time timestamp.
# ---------------------------------------------------------------------
$wsAms.Range("A20").ClearContents()
$wsAms.Range("A20").NumberFormat = "@"
$wsAms.Range("A20").Value = '2021-10-28'

$wsAms.Range("B20").Value = 44497.87484469907

$wsAms.Range("C20").ClearContents()
$wsAms.Range("C20").Value = '152_livetest'

$wsAms.Range("D20").ClearContents()
$wsAms.Range("D20").Value = 155

$wsAms.Range("E20").ClearContents()
$wsAms.Range("E20").Value = 153

$wsAms.Range("F20").ClearContents()
$wsAms.Range("F20").Value = 2

$wsAms.Range("G20").ClearContents()
$wsAms.Range("G20").Value = 4.16

# ---------------------------------------------------------------------
# New rows
# ---------------------------------------------------------------------
    $wsAmsin.Range("A23:A32").NumberFormat = "@"
    $wsAmsin.Range("B23:B32").NumberFormat = "YYYY-MM-DD HH:MM:SS"

    # row 23
    $wsAmsin.Range("A23").Value = '2021-11-17'
    $wsAmsin.Range("B23").Value = 44517.61623751157
    $wsAmsin.Range("C23").Value = '153rsgn'
    $wsAmsin.Range("D23").Value = 155
    $wsAmsin.Range("E23").Value = 137
    $wsAmsin.Range("F23").Value = 18
    $wsAmsin.Range("G23").Value = 9.07

    # row 24
    $wsAmsin.Range("A24").Value = '2021-11-18'
    $wsAmsin.Range("B24").Value = 44518.43382762731
    $wsAmsin.Range("C24").Value = '153fnlrgsn'
    $wsAmsin.Range("D24").Value = 155
    $wsAmsin.Range("E24").Value = 153
    $wsAmsin.Range("F24").Value = 2
    $wsAmsin.Range("G24").Value = 6.86

    # row 25
    $wsAmsin.Range("A25").Value = '2021-12-03'
    $wsAmsin.Range("B25").Value = 44533.72297359954
    $wsAmsin.Range("C25").Value = '154_scndrgsn'
    $wsAmsin.Range("D25").Value = 155
    $wsAmsin.Range("E25").Value = 155
    $wsAmsin.Range("F25").Value = 0
    $wsAmsin.Range("G25").Value = 3.86

    # row 26
    $wsAmsin.Range("A26").Value = '2021-12-06'
    $wsAmsin.Range("B26").Value = 44536.47749277778
    $wsAmsin.Range("C26").Value = '154_fnlrgsn'
    $wsAmsin.Range("D26").Value = 155
    $wsAmsin.Range("E26").Value = 145
    $wsAmsin.Range("F26").Value = 10
    $wsAmsin.Range("G26").Value = 6.99

    # row 27
    $wsAmsin.Range("A27").Value = '2022-01-03'
    $wsAmsin.Range("B27").Value = 44564.76343357639
    $wsAmsin.Range("C27").Value = 'lodash'
    $wsAmsin.Range("D27").Value = 155
    $wsAmsin.Range("E27").Value = 151
    $wsAmsin.Range("F27").Value = 4
    $wsAmsin.Range("G27").Value = 4.57

    # row 28
    $wsAmsin.Range("A28").Value = '2022-01-19'
    $wsAmsin.Range("B28").Value = 44580.76535050926
    $wsAmsin.Range("C28").Value = '165_secondcyc'
    $wsAmsin.Range("D28").Value = 155
    $wsAmsin.Range("E28").Value = 146
    $wsAmsin.Range("F28").Value = 9
    $wsAmsin.Range("G28").Value = 7.04

    # row 29
    $wsAmsin.Range("A29").Value = '2022-01-20'
    $wsAmsin.Range("B29").Value = 44581.40616267361
    $wsAmsin.Range("C29").Value = '156_fnlrsgn'
    $wsAmsin.Range("D29").Value = 155
    $wsAmsin.Range("E29").Value = 153
    $wsAmsin.Range("F29").Value = 2
    $wsAmsin.Range("G29").Value = 4.01

    # row 30
    $wsAmsin.Range("A30").Value = '2022-01-28'
    $wsAmsin.Range("B30").Value = 44589.59209490741
    $wsAmsin.Range("C30").Value = '156audit'
    $wsAmsin.Range("D30").Value = 155
    $wsAmsin.Range("E30").Value = 155
    $wsAmsin.Range("F30").Value = 0
    $wsAmsin.Range("G30").Value = 4.18

    # row 31
    $wsAmsin.Range("A31").Value = '2022-02-07'
    $wsAmsin.Range("B31").Value = 44599.77019583333
    $wsAmsin.Range("C31").Value = 'secondcycle_157'
    $wsAmsin.Range("D31").Value = 155
    $wsAmsin.Range("E31").Value = 153
    $wsAmsin.Range("F31").Value = 2
    $wsAmsin.Range("G31").Value = 4.12

    # row 32
    $wsAmsin.Range("A32").Value = '2022-02-08'
    $wsAmsin.Range("B32").Value = 44600.42689202546
    $wsAmsin.Range("C32").Value = '157_fnl'
    $wsAmsin.Range("D32").Value = 155
    $wsAmsin.Range("E32").Value = 152
    $wsAmsin.Range("F32").Value = 3
    $wsAmsin.Range("G32").Value = 5.53

    $wsBeta.Range("A15:A19").NumberFormat = "@"
    $wsBeta.Range("B15:B19").NumberFormat = "YYYY-MM-DD HH:MM:SS"

    # row 15
    $wsBeta.Range("A15").Value = '2021-11-18'
    $wsBeta.Range("B15").Value = 44518.61075232639
    $wsBeta.Range("C15").Value = '153_beta'
    $wsBeta.Range("D15").Value = 155
    $wsBeta.Range("E15").Value = 150
    $wsBeta.Range("F15").Value = 5
    $wsBeta.Range("G15").Value = 3.96

    # row 16
    $wsBeta.Range("A16").Value = '2021-12-06'
    $wsBeta.Range("B16").Value = 44536.55479989584
    $wsBeta.Range("C16").Value = '154_beta'
    $wsBeta.Range("D16").Value = 155
    $wsBeta.Range("E16").Value = 155
    $wsBeta.Range("F16").Value = 0
    $wsBeta.Range("G16").Value = 4.03

    # row 17
    $wsBeta.Range("A17").Value = '2021-12-23'
    $wsBeta.Range("B17").Value = 44553.52397479166
    $wsBeta.Range("C17").Value = '155_beta'
    $wsBeta.Range("D17").Value = 155
    $wsBeta.Range("E17").Value = 153
    $wsBeta.Range("F17").Value = 2
    $wsBeta.Range("G17").Value = 3.78

    # row 18
    $wsBeta.Range("A18").Value = '2022-01-20'
    $wsBeta.Range("B18").Value = 44581.53989738426
    $wsBeta.Range("C18").Value = '156_beta'
    $wsBeta.Range("D18").Value = 155
    $wsBeta.Range("E18").Value = 154
    $wsBeta.Range("F18").Value = 1
    $wsBeta.Range("G18").Value = 3.41

    # row 19
    $wsBeta.Range("A19").Value = '2022-02-08'
    $wsBeta.Range("B19").Value = 44600.65054744451
    $wsBeta.Range("C19").Value = '157_beta'
    $wsBeta.Range("D19").Value = 155
    $wsBeta.Range("E19").Value = 154
    $wsBeta.Range("F19").Value = 1
    $wsBeta.Range("G19").Value = 3.15

    $wsAms.Range("A21:A26").NumberFormat = "@"
    $wsAms.Range("B21:B26").NumberFormat = "YYYY-MM-DD HH:MM:SS"

    # row 21
    $wsAms.Range("A21").Value = '2021-11-23'
    $wsAms.Range("B21").Value = 44523.46302861111
    $wsAms.Range("C21").Value = '153_live'
    $wsAms.Range("D21").Value = 155
    $wsAms.Range("E21").Value = 155
    $wsAms.Range("F21").Value = 0
    $wsAms.Range("G21").Value = 3.62

    # row 22
    $wsAms.Range("A22").Value = '2021-12-06'
    $wsAms.Range("B22").Value = 44536.87332621527
    $wsAms.Range("C22").Value = '154_live'
    $wsAms.Range("D22").Value = 155
    $wsAms.Range("E22").Value = 153
    $wsAms.Range("F22").Value = 2
    $wsAms.Range("G22").Value = 3.78

    # row 23
    $wsAms.Range("A23").Value = '2021-12-09'
    $wsAms.Range("B23").Value = 44539.78707761574
    $wsAms.Range("C23").Value = '154htfxx'
    $wsAms.Range("D23").Value = 155
    $wsAms.Range("E23").Value = 155
    $wsAms.Range("F23").Value = 0
    $wsAms.Range("G23").Value = 3.43

    # row 24
    $wsAms.Range("A24").Value = '2021-12-23'
    $wsAms.Range("B24").Value = 44553.80614680555
    $wsAms.Range("C24").Value = '155_live'
    $wsAms.Range("D24").Value = 155
    $wsAms.Range("E24").Value = 152
    $wsAms.Range("F24").Value = 3
    $wsAms.Range("G24").Value = 3.79

    # row 25
    $wsAms.Range("A25").Value = '2021-12-27'
    $wsAms.Range("B25").Value = 44557.70875483796
    $wsAms.Range("C25").Value = '155hftfxx'
    $wsAms.Range("D25").Value = 155
    $wsAms.Range("E25").Value = 155
    $wsAms.Range("F25").Value = 0
    $wsAms.Range("G25").Value = 3.92

    # row 26
    $wsAms.Range("A26").Value = '2022-01-20'
    $wsAms.Range("B26").Value = 44581.82316690972
    $wsAms.Range("C26").Value = '156_live'
    $wsAms.Range("D26").Value = 155
    $wsAms.Range("E26").Value = 155
    $wsAms.Range("F26").Value = 0
    $wsAms.Range("G26").Value = 3.44

